# "Generate Report for Archive"
#
# The localization-status report is regenerated: rows whose Status was
# "Ready for handoff" are now "In Translation", and the narrower status
# text lets the Status column(s) shrink accordingly on every sheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $value = $cell.Value()
        # NB: the literal string must be the left-hand operand of -eq;
        # PowerShell otherwise coerces it to the type of a non-string
        # right-hand side (e.g. booleans), causing false matches.
        if ("Ready for handoff" -eq $value) {
            $cell.Value = "In Translation"
        }
    }
}

# Overview sheet: the zh-cn / de-de status columns (E, F) get narrower now
# that the status text is shorter.
# (ColumnWidth is quantized by Excel to whole-pixel steps, so 12.5 here is
# the closest setting that lands the stored column width on the desired
# ~13.4-character target.)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# Per-locale detail sheets: the Status column (C) gets narrower too.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
